$wb = $excel.ActiveWorkbook

# Add the three new sheets at the end of the workbook, in order:
#   同位码 (sheetId 6), h_weima (sheetId 7), l_weima (sheetId 8)
$lastSheet = $wb.Worksheets.Item($wb.Worksheets.Count)

$ws6 = $wb.Worksheets.Add($null, $lastSheet)
$ws6.Name = "同位码"
$ws6.Cells.Item(1, 1).Value = 1
$ws6.Cells.Item(1, 2).Value = 11
$ws6.Cells.Item(2, 1).Value = 1
$ws6.Cells.Item(2, 2).Value = 21
$ws6.Cells.Item(3, 1).Value = 1
$ws6.Cells.Item(3, 2).Value = 31
$ws6.Cells.Item(4, 1).Value = 11
$ws6.Cells.Item(4, 2).Value = 21
$ws6.Cells.Item(5, 1).Value = 11
$ws6.Cells.Item(5, 2).Value = 31
$ws6.Cells.Item(6, 1).Value = 21
$ws6.Cells.Item(6, 2).Value = 31
$ws6.Cells.Item(7, 1).Value = 2
$ws6.Cells.Item(7, 2).Value = 12
$ws6.Cells.Item(8, 1).Value = 2
$ws6.Cells.Item(8, 2).Value = 22
$ws6.Cells.Item(9, 1).Value = 2
$ws6.Cells.Item(9, 2).Value = 32
$ws6.Cells.Item(10, 1).Value = 12
$ws6.Cells.Item(10, 2).Value = 22
$ws6.Cells.Item(11, 1).Value = 12
$ws6.Cells.Item(11, 2).Value = 32
$ws6.Cells.Item(12, 1).Value = 22
$ws6.Cells.Item(12, 2).Value = 32
$ws6.Cells.Item(13, 1).Value = 3
$ws6.Cells.Item(13, 2).Value = 13
$ws6.Cells.Item(14, 1).Value = 3
$ws6.Cells.Item(14, 2).Value = 23
$ws6.Cells.Item(15, 1).Value = 3
$ws6.Cells.Item(15, 2).Value = 33
$ws6.Cells.Item(16, 1).Value = 13
$ws6.Cells.Item(16, 2).Value = 23
$ws6.Cells.Item(17, 1).Value = 13
$ws6.Cells.Item(17, 2).Value = 33
$ws6.Cells.Item(18, 1).Value = 23
$ws6.Cells.Item(18, 2).Value = 33
$ws6.Cells.Item(19, 1).Value = 4
$ws6.Cells.Item(19, 2).Value = 14
$ws6.Cells.Item(20, 1).Value = 4
$ws6.Cells.Item(20, 2).Value = 24
$ws6.Cells.Item(21, 1).Value = 14
$ws6.Cells.Item(21, 2).Value = 24
$ws6.Cells.Item(22, 1).Value = 5
$ws6.Cells.Item(22, 2).Value = 15
$ws6.Cells.Item(23, 1).Value = 5
$ws6.Cells.Item(23, 2).Value = 25
$ws6.Cells.Item(24, 1).Value = 15
$ws6.Cells.Item(24, 2).Value = 23
$ws6.Cells.Item(25, 1).Value = 6
$ws6.Cells.Item(25, 2).Value = 16
$ws6.Cells.Item(26, 1).Value = 6
$ws6.Cells.Item(26, 2).Value = 26
$ws6.Cells.Item(27, 1).Value = 16
$ws6.Cells.Item(27, 2).Value = 26
$ws6.Cells.Item(28, 1).Value = 7
$ws6.Cells.Item(28, 2).Value = 17
$ws6.Cells.Item(29, 1).Value = 7
$ws6.Cells.Item(29, 2).Value = 27
$ws6.Cells.Item(30, 1).Value = 17
$ws6.Cells.Item(30, 2).Value = 27
$ws6.Cells.Item(31, 1).Value = 8
$ws6.Cells.Item(31, 2).Value = 18
$ws6.Cells.Item(32, 1).Value = 8
$ws6.Cells.Item(32, 2).Value = 28
$ws6.Cells.Item(33, 1).Value = 18
$ws6.Cells.Item(33, 2).Value = 28
$ws6.Cells.Item(34, 1).Value = 9
$ws6.Cells.Item(34, 2).Value = 19
$ws6.Cells.Item(35, 1).Value = 9
$ws6.Cells.Item(35, 2).Value = 29
$ws6.Cells.Item(36, 1).Value = 19
$ws6.Cells.Item(36, 2).Value = 29

$ws6.PageSetup.LeftMargin = 54
$ws6.PageSetup.RightMargin = 54
$ws6.PageSetup.TopMargin = 72
$ws6.PageSetup.BottomMargin = 72
$ws6.PageSetup.HeaderMargin = 36
$ws6.PageSetup.FooterMargin = 36

$ws7 = $wb.Worksheets.Add($null, $ws6)
$ws7.Name = "h_weima"

$ws7.PageSetup.LeftMargin = 54
$ws7.PageSetup.RightMargin = 54
$ws7.PageSetup.TopMargin = 72
$ws7.PageSetup.BottomMargin = 72
$ws7.PageSetup.HeaderMargin = 36
$ws7.PageSetup.FooterMargin = 36
$ws7.Cells.Item(1, 1).Value = 5
$ws7.Cells.Item(2, 1).Value = 6
$ws7.Cells.Item(3, 1).Value = 7
$ws7.Cells.Item(4, 1).Value = 8
$ws7.Cells.Item(5, 1).Value = 9
$ws7.Cells.Item(6, 1).Value = 15
$ws7.Cells.Item(7, 1).Value = 16
$ws7.Cells.Item(8, 1).Value = 17
$ws7.Cells.Item(9, 1).Value = 18
$ws7.Cells.Item(10, 1).Value = 19
$ws7.Cells.Item(11, 1).Value = 25
$ws7.Cells.Item(12, 1).Value = 26
$ws7.Cells.Item(13, 1).Value = 27
$ws7.Cells.Item(14, 1).Value = 28
$ws7.Cells.Item(15, 1).Value = 29

$ws8 = $wb.Worksheets.Add($null, $ws7)
$ws8.Name = "l_weima"

$ws8.PageSetup.LeftMargin = 54
$ws8.PageSetup.RightMargin = 54
$ws8.PageSetup.TopMargin = 72
$ws8.PageSetup.BottomMargin = 72
$ws8.PageSetup.HeaderMargin = 36
$ws8.PageSetup.FooterMargin = 36
$ws8.Cells.Item(1, 1).Value = 1
$ws8.Cells.Item(2, 1).Value = 2
$ws8.Cells.Item(3, 1).Value = 3
$ws8.Cells.Item(4, 1).Value = 4
$ws8.Cells.Item(5, 1).Value = 10
$ws8.Cells.Item(6, 1).Value = 11
$ws8.Cells.Item(7, 1).Value = 12
$ws8.Cells.Item(8, 1).Value = 13
$ws8.Cells.Item(9, 1).Value = 14
$ws8.Cells.Item(10, 1).Value = 20
$ws8.Cells.Item(11, 1).Value = 21
$ws8.Cells.Item(12, 1).Value = 22
$ws8.Cells.Item(13, 1).Value = 23
$ws8.Cells.Item(14, 1).Value = 24
$ws8.Cells.Item(15, 1).Value = 30
$ws8.Cells.Item(16, 1).Value = 31
$ws8.Cells.Item(17, 1).Value = 32
$ws8.Cells.Item(18, 1).Value = 33

# Set the active-cell selection for each new sheet; selecting the last
# one (l_weima) last also makes it the workbook's active/selected tab,
# matching the target (activeTab points at l_weima, sheet "六连码" loses
# its previous tabSelected flag automatically).
[void]$ws6.Range("P25").Select()
[void]$ws7.Range("B17").Select()
[void]$ws8.Range("A18").Select()

Write-Host "Added sheets: 同位码, h_weima, l_weima"
